$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.982.48'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.972.17'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.93%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '353.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '112.51'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.69%  '
$ws.Range("E7").Value = '  +0.40%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.631'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.77'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0895'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.71%  '
$ws.Range("E12").Value = '  +0.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.00'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.09%  '
$ws.Range("E14").Value = '  +1.27%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.443.02'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.987.84'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.995'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '52.078.36'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.71'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.54'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '3.32'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.90%  '
$ws.Range("E22").Value = '  +1.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.39'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '270.70'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.81'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.21%  '
$ws.Range("E26").Value = '  +8.89%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '27.76'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.67'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +20.81%  '
$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.111'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +23.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '10.74'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '37.71'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.23'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +10.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '52.88'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.32%  '
$ws.Range("E35").Value = '  -1.02%  '
$ws.Range("E36").Value = '  -4.34%  '
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.45'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.80%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.03'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.31%  '
$ws.Range("E40").Value = '  +1.49%  '
$ws.Range("E41").Value = '  +3.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '23.93'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.73%  '
$ws.Range("E43").Value = '  +1.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.19'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.94%  '
$ws.Range("E45").Value = '  -0.20%  '
$ws.Range("E46").Value = '  +1.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.180.51'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '114.06'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.97%  '
$ws.Range("E49").Value = '  -0.13%  '
$ws.Range("E50").Value = '  +6.32%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.938'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.45%  '
